$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '29.867.42'
$ws.Range("E2").Value = '  -0.46%  '
$ws.Range("D3").Value = '1.895.29'
$ws.Range("E3").Value = '  +0.01%  '
$ws.Range("D4").Value = "'1.0000"
$ws.Range("E4").Value = '  -0.05%  '
$ws.Range("E5").Value = '  -4.88%  '
$ws.Range("D6").Value = "'243.71"
$ws.Range("E6").Value = '  +0.56%  '
$ws.Range("D7").Value = "'0.9996"
$ws.Range("E7").Value = '  -0.05%  '
$ws.Range("D8").Value = "'0.3158"
$ws.Range("E8").Value = '  -3.82%  '
$ws.Range("D9").Value = "'25.40"
$ws.Range("E9").Value = '  -4.42%  '
$ws.Range("D10").Value = "'0.07092"
$ws.Range("E10").Value = '  +0.46%  '
$ws.Range("D11").Value = "'0.08097"
$ws.Range("E11").Value = '  +0.53%  '
$ws.Range("D12").Value = "'0.7697"
$ws.Range("E12").Value = '  +1.86%  '
$ws.Range("B13").Value = 'Polkadot'
$ws.Range("C13").Value = 'https://coinranking.com/coin/25W7FG7om+polkadot-dot'
$ws.Range("D13").Value = "'5.579"
$ws.Range("E13").Value = '  +6.51%  '
$ws.Range("B14").Value = 'WrappedEther'
$ws.Range("C14").Value = 'https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth'
$ws.Range("D14").Value = '1.912.60'
$ws.Range("E14").Value = '  +1.00%  '
$ws.Range("D15").Value = "'92.56"
$ws.Range("E15").Value = '  +0.34%  '
$ws.Range("D16").Value = '29.876.00'
$ws.Range("E16").Value = '  -0.43%  '
$ws.Range("D17").Value = "'6.061"
$ws.Range("E17").Value = '  +2.85%  '
$ws.Range("D18").Value = "'13.91"
$ws.Range("E18").Value = '  -1.30%  '
$ws.Range("D19").Value = "'244.69"
$ws.Range("E19").Value = '  +0.11%  '
$ws.Range("D20").Value = "'0.000007737"
$ws.Range("E20").Value = '  -0.45%  '
$ws.Range("D21").Value = "'8.312"
$ws.Range("E21").Value = '  +19.14%  '
$ws.Range("D22").Value = "'1.000"
$ws.Range("E22").Value = '  +0.14%  '
$ws.Range("D23").Value = '2.131.35'
$ws.Range("E23").Value = '  -0.49%  '
$ws.Range("D24").Value = "'1.000"
$ws.Range("E24").Value = '  -0.03%  '
$ws.Range("D25").Value = "'0.1675"
$ws.Range("E25").Value = '  +3.55%  '
$ws.Range("D26").Value = "'9.434"
$ws.Range("E26").Value = '  +2.18%  '
$ws.Range("D27").Value = "'164.98"
$ws.Range("E27").Value = '  -0.79%  '
$ws.Range("D28").Value = "'18.72"
$ws.Range("E28").Value = '  -0.72%  '
$ws.Range("D29").Value = "'2.062"
$ws.Range("E29").Value = '  -1.65%  '
$ws.Range("D30").Value = "'1.408"
$ws.Range("E30").Value = '  +3.12%  '
$ws.Range("D31").Value = "'1.545"
$ws.Range("E31").Value = '  +1.44%  '
$ws.Range("D32").Value = "'4.538"
$ws.Range("E32").Value = '  +5.78%  '
$ws.Range("D33").Value = "'0.05659"
$ws.Range("E33").Value = '  -4.16%  '
$ws.Range("D34").Value = "'4.071"
$ws.Range("E34").Value = '  -0.20%  '
$ws.Range("D35").Value = "'1.288"
$ws.Range("E35").Value = '  +1.23%  '
$ws.Range("D36").Value = "'0.7425"
$ws.Range("E36").Value = '  +1.72%  '
$ws.Range("D37").Value = "'1.002"
$ws.Range("E37").Value = '  +0.36%  '
$ws.Range("D38").Value = "'2.633"
$ws.Range("E38").Value = '  -3.07%  '
$ws.Range("D39").Value = "'0.01931"
$ws.Range("E39").Value = '  +0.81%  '
$ws.Range("D40").Value = "'2.782"
$ws.Range("E40").Value = '  +0.05%  '
$ws.Range("B41").Value = 'Maker'
$ws.Range("C41").Value = 'https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr'
$ws.Range("D41").Value = '1.125.05'
$ws.Range("E41").Value = '  +13.73%  '
$ws.Range("B42").Value = 'TheSandbox'
$ws.Range("C42").Value = 'https://coinranking.com/coin/pxtKbG5rg+thesandbox-sand'
$ws.Range("D42").Value = "'0.4425"
$ws.Range("E42").Value = '  -0.56%  '
$ws.Range("B43").Value = 'Aave'
$ws.Range("C43").Value = 'https://coinranking.com/coin/ixgUfzmLR+aave-aave'
$ws.Range("D43").Value = "'73.54"
$ws.Range("E43").Value = '  +1.81%  '
$ws.Range("D44").Value = "'0.8546"
$ws.Range("E44").Value = '  +0.25%  '
$ws.Range("D45").Value = "'5.837"
$ws.Range("E45").Value = '  -0.98%  '
$ws.Range("B46").Value = 'Quant'
$ws.Range("C46").Value = 'https://coinranking.com/coin/bauj_21eYVwso+quant-qnt'
$ws.Range("D46").Value = "'104.72"
$ws.Range("E46").Value = '  +3.29%  '
$ws.Range("B47").Value = 'PaxDollar'
$ws.Range("C47").Value = 'https://coinranking.com/coin/JCKLgWPAF+paxdollar-usdp'
$ws.Range("D47").Value = "'0.9992"
$ws.Range("E47").Value = '  -0.05%  '
$ws.Range("B48").Value = 'RenderToken'
$ws.Range("C48").Value = 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
$ws.Range("D48").Value = "'1.884"
$ws.Range("E48").Value = '  +0.41%  '
$ws.Range("B49").Value = 'EnergySwap'
$ws.Range("C49").Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
$ws.Range("D49").Value = "'9.996"
$ws.Range("E49").Value = '  +2.28%  '
$ws.Range("D50").Value = "'7.455"
$ws.Range("E50").Value = '  -1.60%  '
$ws.Range("D51").Value = '2.033.44'
$ws.Range("E51").Value = '  -0.32%  '
